$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every changed cell in the source sheet is stored as literal text (inline
# strings), even where the text happens to look numeric (e.g. "0.999",
# "155.20") or has multiple dots (e.g. "66.323.20"). A plain
# `Range.Value = "..."` lets Excels COM layer auto-convert number-looking
# text into a real Number (dropping significant trailing zeros, switching to
# scientific notation, etc.) and, if coerced back to text via NumberFormat,
# it permanently stamps a new cell style onto every touched cell.
#
# To write a guaranteed-text value without perturbing each cells existing
# style, stage the value (apostrophe-prefixed, forcing text) in an unused
# scratch cell, copy it, and PasteSpecial xlPasteValues (-4163) onto the
# destination: only the literal value/type moves, the destinations own
# number format / style is left untouched.
$scratch = $ws.Range("ZZ1")

$scratch.Value = '''66.323.20'
$scratch.Copy()
$ws.Range('D2').PasteSpecial(-4163)
$scratch.Value = '''  +0.57%  '
$scratch.Copy()
$ws.Range('E2').PasteSpecial(-4163)
$scratch.Value = '''3.179.91'
$scratch.Copy()
$ws.Range('D3').PasteSpecial(-4163)
$scratch.Value = '''  -0.96%  '
$scratch.Copy()
$ws.Range('E3').PasteSpecial(-4163)
$scratch.Value = '''0.999'
$scratch.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$scratch.Value = '''  -0.07%  '
$scratch.Copy()
$ws.Range('E4').PasteSpecial(-4163)
$scratch.Value = '''605.61'
$scratch.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$scratch.Value = '''  +0.45%  '
$scratch.Copy()
$ws.Range('E5').PasteSpecial(-4163)
$scratch.Value = '''155.20'
$scratch.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$scratch.Value = '''  +2.25%  '
$scratch.Copy()
$ws.Range('E6').PasteSpecial(-4163)
$scratch.Value = '''0.999'
$scratch.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$scratch.Value = '''  +0.02%  '
$scratch.Copy()
$ws.Range('E7').PasteSpecial(-4163)
$scratch.Value = '''3.179.15'
$scratch.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$scratch.Value = '''  -0.97%  '
$scratch.Copy()
$ws.Range('E8').PasteSpecial(-4163)
$scratch.Value = '''0.549'
$scratch.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$scratch.Value = '''  +2.07%  '
$scratch.Copy()
$ws.Range('E9').PasteSpecial(-4163)
$scratch.Value = '''0.158'
$scratch.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$scratch.Value = '''  -2.08%  '
$scratch.Copy()
$ws.Range('E10').PasteSpecial(-4163)
$scratch.Value = '''5.71'
$scratch.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$scratch.Value = '''  -7.07%  '
$scratch.Copy()
$ws.Range('E11').PasteSpecial(-4163)
$scratch.Value = '''0.510'
$scratch.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$scratch.Value = '''  -0.19%  '
$scratch.Copy()
$ws.Range('E12').PasteSpecial(-4163)
$scratch.Value = '''0.0000266'
$scratch.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$scratch.Value = '''  -2.15%  '
$scratch.Copy()
$ws.Range('E13').PasteSpecial(-4163)
$scratch.Value = '''38.72'
$scratch.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$scratch.Value = '''  +0.26%  '
$scratch.Copy()
$ws.Range('E14').PasteSpecial(-4163)
$scratch.Value = '''3.697.72'
$scratch.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$scratch.Value = '''  -1.11%  '
$scratch.Copy()
$ws.Range('E15').PasteSpecial(-4163)
$scratch.Value = '''66.367.97'
$scratch.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$scratch.Value = '''  +0.51%  '
$scratch.Copy()
$ws.Range('E16').PasteSpecial(-4163)
$scratch.Value = '''7.41'
$scratch.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$scratch.Value = '''  -0.01%  '
$scratch.Copy()
$ws.Range('E17').PasteSpecial(-4163)
$scratch.Value = '''3.174.79'
$scratch.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$scratch.Value = '''  -1.25%  '
$scratch.Copy()
$ws.Range('E18').PasteSpecial(-4163)
$scratch.Value = '''  +0.20%  '
$scratch.Copy()
$ws.Range('E19').PasteSpecial(-4163)
$scratch.Value = '''512.18'
$scratch.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$scratch.Value = '''  -0.07%  '
$scratch.Copy()
$ws.Range('E20').PasteSpecial(-4163)
$scratch.Value = '''15.52'
$scratch.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$scratch.Value = '''  -1.59%  '
$scratch.Copy()
$ws.Range('E21').PasteSpecial(-4163)
$scratch.Value = '''0.733'
$scratch.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$scratch.Value = '''  -0.61%  '
$scratch.Copy()
$ws.Range('E22').PasteSpecial(-4163)
$scratch.Value = '''8.18'
$scratch.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$scratch.Value = '''  +2.35%  '
$scratch.Copy()
$ws.Range('E23').PasteSpecial(-4163)
$scratch.Value = '''14.93'
$scratch.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$scratch.Value = '''  -1.81%  '
$scratch.Copy()
$ws.Range('E24').PasteSpecial(-4163)
$scratch.Value = '''84.55'
$scratch.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$scratch.Value = '''  -1.00%  '
$scratch.Copy()
$ws.Range('E25').PasteSpecial(-4163)
$scratch.Value = '''0.998'
$scratch.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$scratch.Value = '''  -0.28%  '
$scratch.Copy()
$ws.Range('E26').PasteSpecial(-4163)
$scratch.Value = '''  -0.66%  '
$scratch.Copy()
$ws.Range('E27').PasteSpecial(-4163)
$scratch.Value = '''9.18'
$scratch.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$scratch.Value = '''  -1.13%  '
$scratch.Copy()
$ws.Range('E28').PasteSpecial(-4163)
$scratch.Value = '''2.41'
$scratch.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$scratch.Value = '''  +7.27%  '
$scratch.Copy()
$ws.Range('E29').PasteSpecial(-4163)
$scratch.Value = '''3.11'
$scratch.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$scratch.Value = '''  +7.67%  '
$scratch.Copy()
$ws.Range('E30').PasteSpecial(-4163)
$scratch.Value = '''7.05'
$scratch.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$scratch.Value = '''  +3.55%  '
$scratch.Copy()
$ws.Range('E31').PasteSpecial(-4163)
$scratch.Value = '''28.09'
$scratch.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$scratch.Value = '''  -0.33%  '
$scratch.Copy()
$ws.Range('E32').PasteSpecial(-4163)
$scratch.Value = '''  -0.10%  '
$scratch.Copy()
$ws.Range('E33').PasteSpecial(-4163)
$scratch.Value = '''  -1.52%  '
$scratch.Copy()
$ws.Range('E34').PasteSpecial(-4163)
$scratch.Value = '''6.54'
$scratch.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$scratch.Value = '''  -1.81%  '
$scratch.Copy()
$ws.Range('E35').PasteSpecial(-4163)
$scratch.Value = '''510.09'
$scratch.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$scratch.Value = '''  +4.50%  '
$scratch.Copy()
$ws.Range('E36').PasteSpecial(-4163)
$scratch.Value = '''54.69'
$scratch.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$scratch.Value = '''  -1.47%  '
$scratch.Copy()
$ws.Range('E37').PasteSpecial(-4163)
$scratch.Value = '''  -3.35%  '
$scratch.Copy()
$ws.Range('E38').PasteSpecial(-4163)
$scratch.Value = '''0.0422'
$scratch.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$scratch.Value = '''  -0.37%  '
$scratch.Copy()
$ws.Range('E39').PasteSpecial(-4163)
$scratch.Value = '''  +5.51%  '
$scratch.Copy()
$ws.Range('E40').PasteSpecial(-4163)
$scratch.Value = '''8.86'
$scratch.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$scratch.Value = '''  -0.35%  '
$scratch.Copy()
$ws.Range('E41').PasteSpecial(-4163)
$scratch.Value = '''PEPE'
$scratch.Copy()
$ws.Range('B42').PasteSpecial(-4163)
$scratch.Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$scratch.Copy()
$ws.Range('C42').PasteSpecial(-4163)
$scratch.Value = '''0.0₃0679'
$scratch.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$scratch.Value = '''  +4.67%  '
$scratch.Copy()
$ws.Range('E42').PasteSpecial(-4163)
$scratch.Value = '''TheGraph'
$scratch.Copy()
$ws.Range('B43').PasteSpecial(-4163)
$scratch.Value = '''https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$scratch.Copy()
$ws.Range('C43').PasteSpecial(-4163)
$scratch.Value = '''0.300'
$scratch.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$scratch.Value = '''  +1.90%  '
$scratch.Copy()
$ws.Range('E43').PasteSpecial(-4163)
$scratch.Value = '''2.82'
$scratch.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$scratch.Value = '''  -6.04%  '
$scratch.Copy()
$ws.Range('E44').PasteSpecial(-4163)
$scratch.Value = '''2.43'
$scratch.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$scratch.Value = '''  -1.36%  '
$scratch.Copy()
$ws.Range('E45').PasteSpecial(-4163)
$scratch.Value = '''2.847.84'
$scratch.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$scratch.Value = '''  -6.02%  '
$scratch.Copy()
$ws.Range('E46').PasteSpecial(-4163)
$scratch.Value = '''28.33'
$scratch.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$scratch.Value = '''  -2.87%  '
$scratch.Copy()
$ws.Range('E47').PasteSpecial(-4163)
$scratch.Value = '''  +3.48%  '
$scratch.Copy()
$ws.Range('E48').PasteSpecial(-4163)
$scratch.Value = '''  +0.34%  '
$scratch.Copy()
$ws.Range('E50').PasteSpecial(-4163)
$scratch.Value = '''2.58'
$scratch.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$scratch.Value = '''  +5.42%  '
$scratch.Copy()
$ws.Range('E51').PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false
